$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Bank ID value from PAU0226 -> PAU0228
$ws.Range("M2").Value = "PAU0228"

# Update the EXPLAIN / preparation text to reference the new Bank ID
$ws.Range("F2").Value = "Username : 32070;`nPassword : bni1234;`nRole : 18/19 - Pimpinan Kelompok Investasi/Pengelola Investasi;`nBank ID : PAU0228"

# Update the view state: scrolled/selected cell moved
$ws.Application.ActiveWindow.ScrollColumn = 7
$ws.Range("S2").Select()
